# Student and Course Validation added
# Remove the invalid/duplicate student record in row 6 (Roll 605) from Sheet1.
# Deleting the entire row shifts all subsequent rows (7-16) up by one,
# which naturally renumbers them to rows 6-15 and shrinks the used range
# from A1:E16 down to A1:E15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete row 6 entirely, shifting the rows below it upward.
$ws.Rows.Item(6).Delete()

# Match the resulting active selection recorded in the saved workbook.
$ws.Range("F5").Select()
